# Applies the "ifoCAST full series evaluation" update:
# Each data row (rows 2-16, columns B-K) is shifted one column to the left
# (the value that used to be one column to the right becomes the new value),
# rows 2-6 receive one brand-new trailing value in column K (a newly
# evaluated forecast horizon), and the cell that falls off the right edge of
# each row is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @('-0.07696973751825081','-1.790759194417562','1.03759617918141','-1.19773642068471','0.1535986372329528','-0.467839292583949','0.032721404606882','-0.6351866677392595','-0.2061856925012563','-0.1485141439230462')
    3  = @('-1.208164657383921','1.620190716215051','-0.6151418836510686','0.7361931742665938','0.114755244449692','0.6153159416405229','-0.05259213070561841','0.3764088445323847','0.4340803931105948','0.5580467639488803')
    4  = @('1.524515675405693','-0.7108169244604263','0.6405181334572361','0.01908020364033419','0.5196409008311652','-0.1482671715149762','0.2807338037230269','0.338405352301237','0.4623717231395225','-0.4310464000952693')
    5  = @('-0.4503978874617036','0.9009371704559588','0.279499240639057','0.780059937829888','0.1121518654837466','0.5411528407217497','0.5988243892999598','0.7227907601382453','-0.1706273630965465','0.5995033638472159')
    6  = @('0.06548217429746761','-0.5559557555194342','-0.0553950583286032','-0.7233031306747446','-0.2943021554367415','-0.2366306068585314','-0.1126642360202459','-1.006082359255038','-0.2359516323112753','-0.5018523531907899')
    7  = @('-0.4558763956168127','0.04468430157401831','-0.6232237707721231','-0.19422279553412','-0.1365512469559099','-0.01258487611762438','-0.9060029993524162','-0.1358722724086538','-0.4017729932881683',$null)
    8  = @('0.3984708076537146','-0.2694372646924268','0.1595637105455762','0.2172352591237863','0.3412016299620719','-0.55221649327272','0.2179142336710425','-0.04798648720847212',$null,$null)
    9  = @('-0.4176146551248722','0.0113863201131309','0.06905786869134101','0.1930242395296265','-0.7003938837051653','0.06973684323859711','-0.1961638776409175',$null,$null,$null)
    10 = @('-0.1460596688910241','-0.08838812031281398','0.03557825052547153','-0.8578398727093204','-0.08770914576555788','-0.3536098666450724',$null,$null,$null,$null)
    11 = @('0.4122713362132568','0.5362377070515423','-0.3571804161832495','0.4129503107605129','0.1470495898809984',$null,$null,$null,$null,$null)
    12 = @('0.2558798265160732','-0.6375382967187186','0.1325924302250437','-0.1333082906544708',$null,$null,$null,$null,$null,$null)
    13 = @('-0.465563438873707','0.3045672880700554','0.03866656719054083',$null,$null,$null,$null,$null,$null,$null)
    14 = @('-0.001295251676276088','-0.2671959725557906',$null,$null,$null,$null,$null,$null,$null,$null)
    15 = @('0.3451339801314955',$null,$null,$null,$null,$null,$null,$null,$null,$null)
    16 = @($null,$null,$null,$null,$null,$null,$null,$null,$null,$null)
}

$columns = @('B','C','D','E','F','G','H','I','J','K')

foreach ($row in 2..16) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = $values[$i]
        $cell = $ws.Range("$col$row")
        if ($null -eq $value) {
            $cell.ClearContents()
        } else {
            $cell.Value = [double]$value
        }
    }
}
